# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (column E) and
# "Correspond Handback DateTime" (column H) timestamps for the
# ed63888e-... file row (row 4) on both the "zh-cn" and "de-de" sheets.
# Rows 4 and 5 share identical timestamp text in the source workbook, so
# both rows are updated to keep them in sync, matching the new values.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E4").Value = "2016-03-12 20:17:38"
$zhcn.Range("H4").Value = "2016-03-12 20:17:55"
$zhcn.Range("E5").Value = "2016-03-12 20:17:38"
$zhcn.Range("H5").Value = "2016-03-12 20:17:55"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E4").Value = "2016-03-12 20:17:41"
$dede.Range("H4").Value = "2016-03-12 20:18:00"
$dede.Range("E5").Value = "2016-03-12 20:17:41"
$dede.Range("H5").Value = "2016-03-12 20:18:00"
